$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "11-8=3"
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "70-54=16"
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "82-62=20"
$cell = $tbl.Cell(1, 4)
$cell.Range.Text = "16+2=18"
$cell = $tbl.Cell(1, 5)
$cell.Range.Text = "52-34=18"
$cell = $tbl.Cell(2, 1)
$cell.Range.Text = "49+0=49"
$cell = $tbl.Cell(2, 2)
$cell.Range.Text = "37+34=71"
$cell = $tbl.Cell(2, 3)
$cell.Range.Text = "55+40=95"
$cell = $tbl.Cell(2, 4)
$cell.Range.Text = "15+2=17"
$cell = $tbl.Cell(2, 5)
$cell.Range.Text = "46+24=70"
$cell = $tbl.Cell(3, 1)
$cell.Range.Text = "53-45=8"
$cell = $tbl.Cell(3, 2)
$cell.Range.Text = "36-23=13"
$cell = $tbl.Cell(3, 3)
$cell.Range.Text = "69-27=42"
$cell = $tbl.Cell(3, 4)
$cell.Range.Text = "42-2=40"
$cell = $tbl.Cell(3, 5)
$cell.Range.Text = "63-27=36"
$cell = $tbl.Cell(4, 1)
$cell.Range.Text = "89-47=42"
$cell = $tbl.Cell(4, 2)
$cell.Range.Text = "39-17=22"
$cell = $tbl.Cell(4, 3)
$cell.Range.Text = "94-47=47"
$cell = $tbl.Cell(4, 4)
$cell.Range.Text = "55-1=54"
$cell = $tbl.Cell(4, 5)
$cell.Range.Text = "0+52=52"
$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "44+34=78"
$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "46+25=71"
$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "36+15=51"
$cell = $tbl.Cell(5, 4)
$cell.Range.Text = "36+58=94"
$cell = $tbl.Cell(5, 5)
$cell.Range.Text = "43+44=87"
$cell = $tbl.Cell(6, 1)
$cell.Range.Text = "71-50=21"
$cell = $tbl.Cell(6, 2)
$cell.Range.Text = "66-41=25"
$cell = $tbl.Cell(6, 3)
$cell.Range.Text = "50+4=54"
$cell = $tbl.Cell(6, 4)
$cell.Range.Text = "63-34=29"
$cell = $tbl.Cell(6, 5)
$cell.Range.Text = "63-38=25"
$cell = $tbl.Cell(7, 1)
$cell.Range.Text = "48+3=51"
$cell = $tbl.Cell(7, 2)
$cell.Range.Text = "3+67=70"
$cell = $tbl.Cell(7, 3)
$cell.Range.Text = "24+60=84"
$cell = $tbl.Cell(7, 4)
$cell.Range.Text = "42+4=46"
$cell = $tbl.Cell(7, 5)
$cell.Range.Text = "92-75=17"
$cell = $tbl.Cell(8, 1)
$cell.Range.Text = "83+0=83"
$cell = $tbl.Cell(8, 2)
$cell.Range.Text = "30+66=96"
$cell = $tbl.Cell(8, 3)
$cell.Range.Text = "48-4=44"
$cell = $tbl.Cell(8, 4)
$cell.Range.Text = "58-43=15"
$cell = $tbl.Cell(8, 5)
$cell.Range.Text = "87-31=56"
$cell = $tbl.Cell(9, 1)
$cell.Range.Text = "92-29=63"
$cell = $tbl.Cell(9, 2)
$cell.Range.Text = "37+4=41"
$cell = $tbl.Cell(9, 3)
$cell.Range.Text = "25+21=46"
$cell = $tbl.Cell(9, 4)
$cell.Range.Text = "40+48=88"
$cell = $tbl.Cell(9, 5)
$cell.Range.Text = "78-70=8"
$cell = $tbl.Cell(10, 1)
$cell.Range.Text = "33-9=24"
$cell = $tbl.Cell(10, 2)
$cell.Range.Text = "44+26=70"
$cell = $tbl.Cell(10, 3)
$cell.Range.Text = "13+81=94"
$cell = $tbl.Cell(10, 4)
$cell.Range.Text = "65-25=40"
$cell = $tbl.Cell(10, 5)
$cell.Range.Text = "26+3=29"
$cell = $tbl.Cell(11, 1)
$cell.Range.Text = "14+45=59"
$cell = $tbl.Cell(11, 2)
$cell.Range.Text = "28+1=29"
$cell = $tbl.Cell(11, 3)
$cell.Range.Text = "42+51=93"
$cell = $tbl.Cell(11, 4)
$cell.Range.Text = "10+3=13"
$cell = $tbl.Cell(11, 5)
$cell.Range.Text = "6+40=46"
$cell = $tbl.Cell(12, 1)
$cell.Range.Text = "51+39=90"
$cell = $tbl.Cell(12, 2)
$cell.Range.Text = "59-37=22"
$cell = $tbl.Cell(12, 3)
$cell.Range.Text = "49+46=95"
$cell = $tbl.Cell(12, 4)
$cell.Range.Text = "85-30=55"
$cell = $tbl.Cell(12, 5)
$cell.Range.Text = "89-17=72"
$cell = $tbl.Cell(13, 1)
$cell.Range.Text = "87-65=22"
$cell = $tbl.Cell(13, 2)
$cell.Range.Text = "10+67=77"
$cell = $tbl.Cell(13, 3)
$cell.Range.Text = "56-28=28"
$cell = $tbl.Cell(13, 4)
$cell.Range.Text = "3+73=76"
$cell = $tbl.Cell(13, 5)
$cell.Range.Text = "53+4=57"
$cell = $tbl.Cell(14, 1)
$cell.Range.Text = "42+0=42"
$cell = $tbl.Cell(14, 2)
$cell.Range.Text = "46+3=49"
$cell = $tbl.Cell(14, 3)
$cell.Range.Text = "90-77=13"
$cell = $tbl.Cell(14, 4)
$cell.Range.Text = "78-14=64"
$cell = $tbl.Cell(14, 5)
$cell.Range.Text = "10+6=16"
$cell = $tbl.Cell(15, 1)
$cell.Range.Text = "34-6=28"
$cell = $tbl.Cell(15, 2)
$cell.Range.Text = "54+34=88"
$cell = $tbl.Cell(15, 3)
$cell.Range.Text = "84+2=86"
$cell = $tbl.Cell(15, 4)
$cell.Range.Text = "80+1=81"
$cell = $tbl.Cell(15, 5)
$cell.Range.Text = "14+33=47"
$cell = $tbl.Cell(16, 1)
$cell.Range.Text = "11+3=14"
$cell = $tbl.Cell(16, 2)
$cell.Range.Text = "65+24=89"
$cell = $tbl.Cell(16, 3)
$cell.Range.Text = "8+74=82"
$cell = $tbl.Cell(16, 4)
$cell.Range.Text = "41-38=3"
$cell = $tbl.Cell(16, 5)
$cell.Range.Text = "33-14=19"
$cell = $tbl.Cell(17, 1)
$cell.Range.Text = "62-35=27"
$cell = $tbl.Cell(17, 2)
$cell.Range.Text = "46-9=37"
$cell = $tbl.Cell(17, 3)
$cell.Range.Text = "91+0=91"
$cell = $tbl.Cell(17, 4)
$cell.Range.Text = "77-46=31"
$cell = $tbl.Cell(17, 5)
$cell.Range.Text = "46+43=89"
$cell = $tbl.Cell(18, 1)
$cell.Range.Text = "86-74=12"
$cell = $tbl.Cell(18, 2)
$cell.Range.Text = "75-53=22"
$cell = $tbl.Cell(18, 3)
$cell.Range.Text = "63-7=56"
$cell = $tbl.Cell(18, 4)
$cell.Range.Text = "63-38=25"
$cell = $tbl.Cell(18, 5)
$cell.Range.Text = "36-34=2"
$cell = $tbl.Cell(19, 1)
$cell.Range.Text = "43+17=60"
$cell = $tbl.Cell(19, 2)
$cell.Range.Text = "43-20=23"
$cell = $tbl.Cell(19, 3)
$cell.Range.Text = "61-60=1"
$cell = $tbl.Cell(19, 4)
$cell.Range.Text = "83+8=91"
$cell = $tbl.Cell(19, 5)
$cell.Range.Text = "51-32=19"
$cell = $tbl.Cell(20, 1)
$cell.Range.Text = "88+2=90"
$cell = $tbl.Cell(20, 2)
$cell.Range.Text = "58-29=29"
$cell = $tbl.Cell(20, 3)
$cell.Range.Text = "65+33=98"
$cell = $tbl.Cell(20, 4)
$cell.Range.Text = "90-79=11"
$cell = $tbl.Cell(20, 5)
$cell.Range.Text = "94-15=79"
